# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with
# the latest scraped values (scheduled GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Price "75.626.40" -> "75.668.22"
$ws.Range("D2").Value = "'75.668.22"
$ws.Range("E2").Value = "  +8.70%  "

# Row 3: Price "2.735.20" -> "2.726.49"
$ws.Range("D3").Value = "'2.726.49"
$ws.Range("E3").Value = "  +12.12%  "

$ws.Range("E4").Value = "  -0.04%  "

# Row 5: Price "188.02" -> "187.44"
$ws.Range("D5").Value = "'187.44"
$ws.Range("E5").Value = "  +12.38%  "

# Row 6: Price "592.71" -> "592.31"
$ws.Range("D6").Value = "'592.31"
$ws.Range("E6").Value = "  +4.68%  "

$ws.Range("E7").Value = "  -0.10%  "

# Row 8: Price "0.544" -> "0.543"
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  +5.45%  "

$ws.Range("E9").Value = "  +15.16%  "

# Row 10: Price "2.734.07" -> "2.728.58"
$ws.Range("D10").Value = "'2.728.58"
$ws.Range("E10").Value = "  +12.27%  "

# Row 11: Price "0.163" -> "0.162"
$ws.Range("D11").Value = "'0.162"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").Value = "  +8.41%  "

# Row 13: Price "4.80" -> "4.79"
$ws.Range("D13").Value = "'4.79"
$ws.Range("E13").Value = "  +1.90%  "

# Row 14: Price "3.183.36" -> "3.232.18"
$ws.Range("D14").Value = "'3.232.18"
$ws.Range("E14").Value = "  +12.25%  "

# Row 15: Price "75.473.59" -> "75.500.13"
$ws.Range("D15").Value = "'75.500.13"
$ws.Range("E15").Value = "  +8.67%  "

$ws.Range("E16").Value = "  +7.22%  "

# Row 17: Price "27.15" -> "27.17"
$ws.Range("D17").Value = "'27.17"
$ws.Range("E17").Value = "  +13.07%  "

# Row 18: Price "2.712.71" -> "2.722.36"
$ws.Range("D18").Value = "'2.722.36"
$ws.Range("E18").Value = "  +11.70%  "

# Row 19: Price "9.54" -> "9.45"
$ws.Range("D19").Value = "'9.45"
$ws.Range("E19").Value = "  +31.84%  "

# Row 20: Price "12.26" -> "12.22"
$ws.Range("D20").Value = "'12.22"
$ws.Range("E20").Value = "  +12.47%  "

# Row 21: Price "378.53" -> "377.89"
$ws.Range("D21").Value = "'377.89"
$ws.Range("E21").Value = "  +9.83%  "

# Row 22: Price "2.32" -> "2.31"
$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +16.37%  "

# Row 23: Price "4.11" -> "4.12"
$ws.Range("D23").Value = "'4.12"
$ws.Range("E23").Value = "  +6.28%  "

$ws.Range("E24").Value = "  +4.43%  "

# Row 25: Price "71.34" -> "71.33"
$ws.Range("D25").Value = "'71.33"
$ws.Range("E25").Value = "  +7.92%  "

$ws.Range("E26").Value = "  -0.02%  "

# Row 27: Price "4.27" -> "4.25"
$ws.Range("D27").Value = "'4.25"
$ws.Range("E27").Value = "  +11.12%  "

# Row 28: Price "9.67" -> "9.68"
$ws.Range("D28").Value = "'9.68"
$ws.Range("E28").Value = "  +14.08%  "

# Row 29: Price "2.845.44" -> "2.858.37"
$ws.Range("D29").Value = "'2.858.37"
$ws.Range("E29").Value = "  +11.91%  "

# Row 30: Price "1.00" -> "0.995"
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -1.23%  "

# Row 31: Price "0.0₃0992" -> "0.0₃0998"
$ws.Range("D31").Value = "'0.0₃0998"
$ws.Range("E31").Value = "  +16.81%  "

# Row 32: Price "526.99" -> "525.88"
$ws.Range("D32").Value = "'525.88"
$ws.Range("E32").Value = "  +15.41%  "

$ws.Range("E33").Value = "  +13.74%  "

# Row 34: Price "7.92" -> "7.91"
$ws.Range("D34").Value = "'7.91"
$ws.Range("E34").Value = "  +7.10%  "

# Row 35: Price "1.79" -> "1.80"
$ws.Range("D35").Value = "'1.80"
$ws.Range("E35").Value = "  +11.39%  "

$ws.Range("E36").Value = "  -0.06%  "

# Row 37: Price "0.120" -> "0.121"
$ws.Range("D37").Value = "'0.121"
$ws.Range("E37").Value = "  +7.66%  "

# Row 38: Price "161.39" -> "161.51"
$ws.Range("D38").Value = "'161.51"
$ws.Range("E38").Value = "  +1.41%  "

# Row 39: Price "19.62" -> "19.61"
$ws.Range("D39").Value = "'19.61"
$ws.Range("E39").Value = "  +7.43%  "

$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("E41").Value = "  -0.01%  "

# Row 42: Price "174.11" -> "174.20"
$ws.Range("D42").Value = "'174.20"
$ws.Range("E42").Value = "  +27.85%  "

# Row 43: Price "5.10" -> "5.09"
$ws.Range("D43").Value = "'5.09"
$ws.Range("E43").Value = "  +15.18%  "

$ws.Range("E44").Value = "  +13.67%  "

$ws.Range("E45").Value = "  +10.30%  "

# Row 46: Price "1.23" -> "1.24"
$ws.Range("D46").Value = "'1.24"
$ws.Range("E46").Value = "  +14.15%  "

# Row 47: Price "2.43" -> "2.42"
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +15.42%  "

# Row 48: Price "39.16" -> "39.18"
$ws.Range("D48").Value = "'39.18"
$ws.Range("E48").Value = "  +3.47%  "

# Row 49: Price "0.0853" -> "0.0857"
$ws.Range("D49").Value = "'0.0857"
$ws.Range("E49").Value = "  +18.81%  "

$ws.Range("E51").Value = "  +12.64%  "
